$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.058.97'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '2.925.63'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '593.83'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.506'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.87'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.144'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.440'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000226'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '33.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = '3.413.30'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '61.180.43'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.71'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '2.922.91'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '431.55'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '81.41'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.03'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  +6.30%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '26.49'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').Value = '0.0₃0851'
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.01'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.64'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('E37').Value = '  +3.96%  '
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.60'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.47'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '375.56'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.11%  '
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').Value = '2.723.37'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '130.67'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '24.04'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.02'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('E51').Value = '  +3.01%  '
